$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume(1h) (E) columns stay as plain text,
# matching the workbook's existing inline-string storage, so Excel
# does not reinterpret values like "26.015.09" or "0.9981" as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.015.09'
$ws.Cells.Item(2, 5).Value = '  -2.94%  '
$ws.Cells.Item(3, 4).Value = '1.839.19'
$ws.Cells.Item(3, 5).Value = '  -1.57%  '
$ws.Cells.Item(4, 4).Value = '0.9981'
$ws.Cells.Item(4, 5).Value = '  -0.27%  '
$ws.Cells.Item(5, 4).Value = '279.04'
$ws.Cells.Item(5, 5).Value = '  -4.75%  '
$ws.Cells.Item(6, 4).Value = '0.9984'
$ws.Cells.Item(6, 5).Value = '  -0.56%  '
$ws.Cells.Item(7, 5).Value = '  -4.01%  '
$ws.Cells.Item(8, 4).Value = '0.3500'
$ws.Cells.Item(8, 5).Value = '  -6.05%  '
$ws.Cells.Item(9, 4).Value = '44.93'
$ws.Cells.Item(9, 5).Value = '  -1.08%  '
$ws.Cells.Item(10, 4).Value = '0.06813'
$ws.Cells.Item(10, 5).Value = '  -4.58%  '
$ws.Cells.Item(11, 4).Value = '19.99'
$ws.Cells.Item(11, 5).Value = '  -6.39%  '
$ws.Cells.Item(12, 4).Value = '0.8069'
$ws.Cells.Item(12, 5).Value = '  -8.83%  '
$ws.Cells.Item(13, 4).Value = '0.07805'
$ws.Cells.Item(13, 5).Value = '  -4.03%  '
$ws.Cells.Item(14, 4).Value = '1.835.75'
$ws.Cells.Item(14, 5).Value = '  -1.47%  '
$ws.Cells.Item(15, 4).Value = '5.079'
$ws.Cells.Item(15, 5).Value = '  -3.75%  '
$ws.Cells.Item(16, 4).Value = '88.44'
$ws.Cells.Item(16, 5).Value = '  -4.03%  '
$ws.Cells.Item(17, 4).Value = '0.9976'
$ws.Cells.Item(17, 5).Value = '  -0.30%  '
$ws.Cells.Item(18, 4).Value = '14.17'
$ws.Cells.Item(18, 5).Value = '  -4.19%  '
$ws.Cells.Item(19, 4).Value = '0.000008078'
$ws.Cells.Item(19, 5).Value = '  -4.64%  '
$ws.Cells.Item(20, 4).Value = '0.9988'
$ws.Cells.Item(20, 5).Value = '  -0.34%  '
$ws.Cells.Item(21, 4).Value = '26.037.52'
$ws.Cells.Item(21, 5).Value = '  -2.68%  '
$ws.Cells.Item(22, 4).Value = '4.771'
$ws.Cells.Item(22, 5).Value = '  -3.73%  '
$ws.Cells.Item(23, 4).Value = '10.06'
$ws.Cells.Item(23, 5).Value = '  -5.34%  '
$ws.Cells.Item(24, 4).Value = '6.214'
$ws.Cells.Item(24, 5).Value = '  -2.15%  '
$ws.Cells.Item(25, 4).Value = '2.361'
$ws.Cells.Item(25, 5).Value = '  +3.60%  '
$ws.Cells.Item(26, 4).Value = '143.12'
$ws.Cells.Item(26, 5).Value = '  -1.86%  '
$ws.Cells.Item(27, 5).Value = '  -4.09%  '
$ws.Cells.Item(28, 4).Value = '17.21'
$ws.Cells.Item(28, 5).Value = '  -4.07%  '
$ws.Cells.Item(29, 4).Value = '109.74'
$ws.Cells.Item(29, 5).Value = '  -3.01%  '
$ws.Cells.Item(30, 4).Value = '4.367'
$ws.Cells.Item(30, 5).Value = '  -6.79%  '
$ws.Cells.Item(31, 4).Value = '4.290'
$ws.Cells.Item(31, 5).Value = '  -6.89%  '
$ws.Cells.Item(32, 4).Value = '0.08805'
$ws.Cells.Item(32, 5).Value = '  -3.10%  '
$ws.Cells.Item(33, 4).Value = '0.04867'
$ws.Cells.Item(33, 5).Value = '  -2.94%  '
$ws.Cells.Item(34, 4).Value = '1.166'
$ws.Cells.Item(34, 5).Value = '  -0.12%  '
$ws.Cells.Item(35, 4).Value = '0.7322'
$ws.Cells.Item(35, 5).Value = '  -8.77%  '
$ws.Cells.Item(36, 4).Value = '2.868'
$ws.Cells.Item(36, 5).Value = '  -3.00%  '
$ws.Cells.Item(37, 4).Value = '3.207'
$ws.Cells.Item(37, 5).Value = '  +0.56%  '
$ws.Cells.Item(38, 4).Value = '0.9986'
$ws.Cells.Item(38, 5).Value = '  -0.23%  '
$ws.Cells.Item(39, 4).Value = '2.392'
$ws.Cells.Item(39, 5).Value = '  -9.48%  '
$ws.Cells.Item(40, 4).Value = '0.01849'
$ws.Cells.Item(40, 5).Value = '  -4.35%  '
$ws.Cells.Item(41, 4).Value = '0.5170'
$ws.Cells.Item(41, 5).Value = '  -13.82%  '
$ws.Cells.Item(42, 4).Value = '0.9469'
$ws.Cells.Item(42, 5).Value = '  -10.83%  '
$ws.Cells.Item(43, 4).Value = '117.34'
$ws.Cells.Item(43, 5).Value = '  +1.53%  '
$ws.Cells.Item(44, 4).Value = '6.264'
$ws.Cells.Item(44, 5).Value = '  -3.18%  '
$ws.Cells.Item(45, 4).Value = '8.013'
$ws.Cells.Item(45, 5).Value = '  -7.99%  '
$ws.Cells.Item(46, 4).Value = '0.9979'
$ws.Cells.Item(46, 5).Value = '  -0.75%  '
$ws.Cells.Item(47, 4).Value = '0.4522'
$ws.Cells.Item(47, 5).Value = '  -13.64%  '
$ws.Cells.Item(48, 4).Value = '0.1364'
$ws.Cells.Item(48, 5).Value = '  -8.10%  '
$ws.Cells.Item(49, 4).Value = '9.305'
$ws.Cells.Item(50, 4).Value = '36.17'
$ws.Cells.Item(50, 5).Value = '  -2.66%  '
$ws.Cells.Item(51, 4).Value = '0.05925'
$ws.Cells.Item(51, 5).Value = '  -2.01%  '
